# Atualização automática dos dados do dashboard (aba "Entrada")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrada")

# Linha 2 - DEVOLUÇÃO
$ws.Range("B2").Value = "R$ 777.126,11"
$ws.Range("D2").Value = "R$ 777.126,11"
$ws.Range("E2").Value = "R$ 777.126,11"

# Linha 3 - FERRAMENTAS/ MATRIZARIA
$ws.Range("B3").Value = "R$ 528.494,31"
$ws.Range("D3").Value = "R$ 528.494,31"
$ws.Range("F3").Value = "77,72 %"

# Linha 5 - REFUGO REAL (PROCESSO)
$ws.Range("B5").Value = "R$ 328.104,17"
$ws.Range("D5").Value = "R$ 328.104,17"
$ws.Range("E5").Value = "R$ 328.104,17"

# Linha 6 - CUSTO DESENVOLVIMENTO
$ws.Range("B6").Value = "R$ 313.457,68"
$ws.Range("D6").Value = "R$ 313.457,68"
$ws.Range("E6").Value = "R$ 313.457,68"

# Linha 7 - agora FRETES (era MANUTENCAO)
$ws.Range("A7").Value = "FRETES"
$ws.Range("B7").Value = "R$ 214.614,62"
$ws.Range("C7").Value = "R$ 0,00"
$ws.Range("D7").Value = "R$ 214.614,62"
$ws.Range("E7").Value = "R$ 376.000,00"
$ws.Range("F7").Value = "57,08 %"

# Linha 8 - agora MANUTENCAO (era FRETES)
$ws.Range("A8").Value = "MANUTENCAO"
$ws.Range("B8").Value = "R$ 206.862,29"
$ws.Range("C8").Value = "R$ 191.636,71"
$ws.Range("D8").Value = "R$ 398.499,00"
$ws.Range("E8").Value = "R$ 480.000,00"
$ws.Range("F8").Value = "83,02 %"

# Linha 9 - REFUGO MP+CP*
$ws.Range("B9").Value = "R$ 195.167,49"
$ws.Range("D9").Value = "R$ 195.167,49"
$ws.Range("F9").Value = "69,70 %"

# Linha 10 - agora DESP. INDUSTRIAL (era OLEOS E LUBRIFICANTES)
$ws.Range("A10").Value = "DESP. INDUSTRIAL"
$ws.Range("B10").Value = "R$ 187.496,71"
$ws.Range("C10").Value = "R$ 109.131,35"
$ws.Range("D10").Value = "R$ 296.628,06"
$ws.Range("E10").Value = "R$ 450.000,00"
$ws.Range("F10").Value = "65,92 %"

# Linha 11 - agora OLEOS E LUBRIFICANTES (era DESP. INDUSTRIAL)
$ws.Range("A11").Value = "OLEOS E LUBRIFICANTES"
$ws.Range("B11").Value = "R$ 82.091,98"
$ws.Range("C11").Value = "R$ 108.767,66"
$ws.Range("D11").Value = "R$ 190.859,64"
$ws.Range("E11").Value = "R$ 280.000,00"
$ws.Range("F11").Value = "68,16 %"

# Linha 12 - EMBALAGENS
$ws.Range("C12").Value = "R$ 97.082,75"
$ws.Range("D12").Value = "R$ 155.647,35"
$ws.Range("F12").Value = "97,28 %"

# Linha 15 - CUSTO COM DESENVOLVIMENTO
$ws.Range("B15").Value = "R$ 8.615,99"
$ws.Range("D15").Value = "R$ 8.615,99"
$ws.Range("E15").Value = "R$ 8.615,99"

# Linha 18 - Total Geral
$ws.Range("B18").Value = "R$ 3.128.596,78"
$ws.Range("C18").Value = "R$ 511.151,23"
$ws.Range("D18").Value = "R$ 3.639.748,01"
$ws.Range("E18").Value = "R$ 4.497.304,95"
$ws.Range("F18").Value = "80,93 %"
